# Update "想去人数" (attendance count) figures that changed between two
# scrapes of the source data.
#
# Sheet "展览" (exhibitions):
#   Row 4 (南昌·鹃歌袂 代号鹃同人only): F4  171 -> 172
#   Row 6 (南昌·CM04动漫游戏博览会):    F6 5488 -> 5491
#
# Sheet "全部类型" (all types, aggregates the same events):
#   Row 4 (南昌·鹃歌袂 代号鹃同人only): F4  171 -> 172
#   Row 7 (南昌·CM04动漫游戏博览会):    F7 5488 -> 5491

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 172
$wsExpo.Range("F6").Value = 5491

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 172
$wsAll.Range("F7").Value = 5491
